$d = $word.ActiveDocument

# The placeholder/ID paragraph is the very first paragraph in the body.
$p1 = $d.Paragraphs(1)

# Add a paragraph border whose lines are not switched on, only the
# "distance from text" (space) is set on all four sides -> <w:pBdr> with
# just w:space="5" on top/left/bottom/right.
$p1.Borders.DistanceFromTop = 5
$p1.Borders.DistanceFromLeft = 5
$p1.Borders.DistanceFromBottom = 5
$p1.Borders.DistanceFromRight = 5

# Left indent goes from 120 twips (6pt) to 225 twips (11.25pt). Word's
# object model works in points, so divide twips by 20.
$p1.LeftIndent = 11.25

# Replace the old placeholder id text (plus its trailing space, which
# lived in a second run) with the new id text, collapsing both runs
# into a single run with no trailing space.
$d.Content.Find.Execute("**ID__AFFARS_mp_5315_3_topic_14__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_MP_5315_3_3_2__ID**", 2)
